$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Packet")

# Insert a new row before row 8 ("Flag") to make room for the new "SkillLv" entry.
$ws.Rows.Item(8).Insert()

# Populate the newly inserted row 8 with the new model field.
$ws.Range("A8").Value = "SkillLv"
$ws.Range("B8").Value = "INT"

# Reflect the author's last selection in the sheet.
$ws.Range("B7").Select()
